$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17:18").Insert()
$ws.Range("A16:H16").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Range("A16:H16").Copy()
$ws.Range("A18:H18").PasteSpecial(-4122)
$ws.Rows("17:18").RowHeight = 36.75

$ws.Range("A17").Value = "JP3"
$ws.Range("A18").Value = "JP3 Conn"
$ws.Range("B17").Value = "SAM8980-ND"
$ws.Range("C17").Value = "Connector Header Surface Mount 2 position 0.100"" (2.54mm)"
$ws.Range("H17").Value = "https://www.digikey.com/en/products/detail/samtec-inc/TSM-102-01-T-SV/2685536"
$ws.Range("B18").Value = "ED3082-ND"
$ws.Range("C18").Value = "2 (1 x 2) Position Shunt Connector Black Closed Top 0.100"" (2.54mm) Gold"
$ws.Range("H18").Value = "https://www.digikey.com/en/products/detail/on-shore-technology-inc/EDJ1G0/2752411"

$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0.34
$ws.Range("F17").Value = "Digikey"
$ws.Range("G17").Formula = "=D17*E17"

$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0.1
$ws.Range("F18").Value = "Digikey"
$ws.Range("G18").Formula = "=D18*E18"
